$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 3, shifting existing rows 3-4 down to 4-5.
$ws.Rows("3").Insert()

# Populate the newly inserted row 3 with "EWAIT" in column A (matching row 5's style/content).
$ws.Range("A3").Value = "EWAIT"

# Copy the style (borders etc.) from the row below into the new row's cells.
$ws.Range("A4:D4").Copy()
$ws.Range("A3:D3").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Update selection to match the target state (single active cell A5).
$ws.Range("A5").Select()
